# "uji robust fix cropping" - update the cropping-percentage column headers in
# row 3 of Sheet1 to describe the crop side/position explicitly, and widen the
# affected columns so the longer labels fit. Also add a short explanatory note
# ("penjelasan uji") and refresh the view (selection/zoom) like a user would
# leave it after finishing the edit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3 header text: old generic "25%/50% Horizontal/Vertikal" labels are
#     replaced with labels that say exactly where the crop was taken from.
$ws.Range("D3").Value = "50% Horizontal Atas"
$ws.Range("E3").Value = "50% Horizontal Bawah"
$ws.Range("F3").Value = "50% Vertikal Kiri"
$ws.Range("G3").Value = "50% Vertikal Kanan"

# --- Column widths: column C keeps its original width, columns D:G (which now
#     hold the longer labels above) are widened so the text isn't truncated.
$ws.Columns("D:G").ColumnWidth = 19.17

# --- View bookkeeping: leave the sheet scrolled/selected/zoomed the way the
#     author left it after finishing the edit.
$excel.ActiveWindow.Zoom = 100
$ws.Range("H7").Select()
